$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "name=login"
$ws.Range("C3").Value = "enter url"

$ws.Range("C3").Select()
